$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> label value to fill into column E
$values = @{
    34 = -1
    35 = 0
    36 = 0
    37 = -1
    38 = -1
    39 = -1
    40 = -1
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = -1
    46 = -1
    47 = -1
    48 = -1
    49 = 0
    50 = -1
    51 = 0
    52 = 0
    53 = -1
    54 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}

# Update view state: scroll position and selection to match edit
$ws.Range("E55").Select()
$excel.ActiveWindow.ScrollRow = 54
